$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new "regular giant" luminosity sub-classes into column D,
# pushing the existing "III (regular giants)" / "IV (sub-giants)" entries down.
$ws.Range("D3").Value = "IIIa (luminous regular giant)"
$ws.Range("D4").Value = "III (regular giants)"
$ws.Range("D5").Value = "IIIb (less luminous regular giant"
$ws.Range("D6").Value = "IV (sub-giants)"

# Drop the (invisible, empty) border formatting that had been applied to the
# header row and the data block, consolidating the style palette down to the
# two still-used cell formats.
$ws.Range("C1:F1").Borders.LineStyle = -4142
$ws.Range("C2:F3").Borders.LineStyle = -4142
$ws.Range("C4:E6").Borders.LineStyle = -4142

# Column D needs to widen to fit the new, longer text (best-fit width).
$ws.Columns("D").ColumnWidth = 28.42

# Restore the cursor / selection position.
$ws.Range("D6").Select()
